# Updates cryptos list figures (price + 1h volume change) to match the latest
# coinranking.com snapshot. Row 30/31 (Toncoin / InjectiveProtocol) swap rank
# order, so their whole rows (name/link/price/volume) are replaced.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'48.374.91"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +2.07%  '
$ws.Range("D3").Value = "'2.522.23"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.46%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").Value = "'323.53"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.25%  '
$ws.Range("D6").Value = "'109.47"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.88%  '
$ws.Range("E7").Value = '  -0.18%  '
$ws.Range("D8").Value = "'0.999"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.06%  '
$ws.Range("E9").Value = '  +3.70%  '
$ws.Range("D10").Value = "'40.33"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.27%  '
$ws.Range("E11").Value = '  +0.37%  '
$ws.Range("D12").Value = "'19.64"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +5.35%  '
$ws.Range("E13").Value = '  +0.79%  '
$ws.Range("E14").Value = '  -0.23%  '
$ws.Range("D15").Value = "'2.908.00"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.17%  '
$ws.Range("D16").Value = "'2.516.61"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.19%  '
$ws.Range("E17").Value = '  -0.75%  '
$ws.Range("D18").Value = "'48.220.42"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.82%  '
$ws.Range("D19").Value = "'13.44"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.69%  '
$ws.Range("E20").Value = '  -0.88%  '
$ws.Range("D21").Value = "'0.0₃0946"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.27%  '
$ws.Range("E22").Value = '  +3.02%  '
$ws.Range("D23").Value = "'72.48"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.27%  '
$ws.Range("D24").Value = "'267.62"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +6.96%  '
$ws.Range("D25").Value = "'2.58"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.30%  '
$ws.Range("D26").Value = "'26.23"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.04%  '
$ws.Range("E27").Value = '  +0.11%  '
$ws.Range("E28").Value = '  +1.62%  '
$ws.Range("D29").Value = "'0.144"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +5.04%  '
$ws.Range("B30").Value = 'InjectiveProtocol'
$ws.Range("C30").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D30").Value = "'35.37"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.79%  '
$ws.Range("B31").Value = 'Toncoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D31").Value = "'2.08"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -9.23%  '
$ws.Range("D32").Value = "'49.97"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.61%  '
$ws.Range("D33").Value = "'20.05"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.44%  '
$ws.Range("E34").Value = '  -0.94%  '
$ws.Range("D35").Value = "'1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.04%  '
$ws.Range("D36").Value = "'0.0789"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.27%  '
$ws.Range("E37").Value = '  -1.01%  '
$ws.Range("D38").Value = "'4.71"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.65%  '
$ws.Range("E39").Value = '  -1.24%  '
$ws.Range("D40").Value = "'0.113"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.15%  '
$ws.Range("D41").Value = "'22.26"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +3.64%  '
$ws.Range("D42").Value = "'118.82"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.18%  '
$ws.Range("E43").Value = '  -3.43%  '
$ws.Range("E44").Value = '  +0.35%  '
$ws.Range("D45").Value = "'2.002.04"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.03%  '
$ws.Range("E46").Value = '  +0.13%  '
$ws.Range("E47").Value = '  -2.83%  '
$ws.Range("E48").Value = '  +3.87%  '
$ws.Range("E49").Value = '  +0.20%  '
$ws.Range("D50").Value = "'5.26"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.51%  '
$ws.Range("D51").Value = "'80.47"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.28%  '
